$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 601.375
$ws.Range("I33").Value = 402.54544
$ws.Range("J33").Value = 1038.8
$ws.Range("K33").Value = 402.54544
$ws.Range("L33").Value = 1038.8
$ws.Range("M33").Value = -173.54544
$ws.Range("N33").Value = -1496.8
$ws.Range("H43").Value = 4000.5
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4138
$ws.Range("H98").Value = 2007.5385
$ws.Range("I98").Value = 1814.7778
$ws.Range("J98").Value = 2441.25
$ws.Range("K98").Value = 1814.7778
$ws.Range("L98").Value = 2441.25
$ws.Range("M98").Value = -316.7778000000001
$ws.Range("N98").Value = -5437.25
$ws.Range("H116").Value = 7848.4
$ws.Range("I116").Value = 5496.25
$ws.Range("K116").Value = 5496.25
$ws.Range("M116").Value = -2054.25
$ws.Range("H122").Value = 2007.5385
$ws.Range("I122").Value = 1814.7778
$ws.Range("J122").Value = 2441.25
$ws.Range("K122").Value = 5444.3334
$ws.Range("L122").Value = 7323.75
$ws.Range("M122").Value = -2994.3334
$ws.Range("N122").Value = -12223.75
$ws.Range("H137").Value = 2198.7585
$ws.Range("I137").Value = 1833
$ws.Range("J137").Value = 4190.1113
$ws.Range("K137").Value = 5499
$ws.Range("L137").Value = 12570.3339
$ws.Range("M137").Value = -2949
$ws.Range("N137").Value = -17670.3339

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2338.2727
$ws.Range("I61").Value = 2338.2727
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2338.2727
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2126.2727
$ws.Range("H136").Value = 2338.2727
$ws.Range("I136").Value = 2338.2727
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7014.8181
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4464.8181
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7883.08
$ws.Range("J107").Value = 16505
$ws.Range("L107").Value = 16505
$ws.Range("N107").Value = -20345

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2518.5
$ws.Range("I132").Value = 2072.3667
$ws.Range("K132").Value = 6217.1001
$ws.Range("M132").Value = -3687.1001
$ws.Range("H134").Value = 10665.667
$ws.Range("J134").Value = 12006
$ws.Range("L134").Value = 36018
$ws.Range("N134").Value = -41088
$ws.Range("H141").Value = 40214.285
$ws.Range("J141").Value = 41000
$ws.Range("L141").Value = 41000
$ws.Range("N141").Value = -51360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1538.8235
$ws.Range("I5").Value = 1490.6666
$ws.Range("J5").Value = 1900
$ws.Range("K5").Value = 4471.9998
$ws.Range("L5").Value = 5700
$ws.Range("M5").Value = -4359.9998
$ws.Range("N5").Value = -5924
$ws.Range("H56").Value = 7290.75
$ws.Range("I56").Value = 7290.75
$ws.Range("K56").Value = 7290.75
$ws.Range("M56").Value = -6760.75
$ws.Range("H98").Value = 1355.75
$ws.Range("I98").Value = 1674
$ws.Range("J98").Value = 401
$ws.Range("K98").Value = 5022
$ws.Range("L98").Value = 1203
$ws.Range("M98").Value = -3524
$ws.Range("N98").Value = -4199
$ws.Range("H135").Value = 1538.8235
$ws.Range("I135").Value = 1490.6666
$ws.Range("J135").Value = 1900
$ws.Range("K135").Value = 13415.9994
$ws.Range("L135").Value = 17100
$ws.Range("M135").Value = -10880.9994
$ws.Range("N135").Value = -22170
$ws.Range("H138").Value = 5083.294
$ws.Range("I138").Value = 4522.9287
$ws.Range("K138").Value = 13568.7861
$ws.Range("M138").Value = -8428.786100000001
$ws.Range("H141").Value = 100002136
$ws.Range("I141").Value = 100002136
$ws.Range("K141").Value = 300006408
$ws.Range("M141").Value = -300001228

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 565
$ws.Range("I19").Value = 672.5
$ws.Range("J19").Value = 350
$ws.Range("K19").Value = 672.5
$ws.Range("L19").Value = 350
$ws.Range("M19").Value = -384.5
$ws.Range("N19").Value = -926
$ws.Range("H126").Value = 8429.666999999999
$ws.Range("I126").Value = 5500
$ws.Range("K126").Value = 16500
$ws.Range("M126").Value = -14030

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2998.139
$ws.Range("J93").Value = 5078.091
$ws.Range("L93").Value = 5078.091
$ws.Range("N93").Value = -7574.091
$ws.Range("H128").Value = 74049.5
$ws.Range("J128").Value = 74049.5
$ws.Range("L128").Value = 74049.5
$ws.Range("N128").Value = -84009.5
$ws.Range("H132").Value = 2487.8809
$ws.Range("I132").Value = 2117.8386
$ws.Range("K132").Value = 6353.5158
$ws.Range("M132").Value = -3823.5158
$ws.Range("H141").Value = 84400
$ws.Range("J141").Value = 78800
$ws.Range("L141").Value = 78800
$ws.Range("N141").Value = -89160

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 532500
$ws.Range("I15").Value = 532500
$ws.Range("K15").Value = 532500
$ws.Range("M15").Value = -532212
$ws.Range("H45").Value = 13775.4
$ws.Range("I45").Value = 7912.3335
$ws.Range("J45").Value = 16288.143
$ws.Range("K45").Value = 7912.3335
$ws.Range("L45").Value = 16288.143
$ws.Range("M45").Value = -7421.3335
$ws.Range("N45").Value = -17270.143
$ws.Range("H126").Value = 2364.4443
$ws.Range("I126").Value = 2410
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 7230
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -4760
$ws.Range("N126").Value = -10940
$ws.Range("H128").Value = 57499.25
$ws.Range("J128").Value = 57499.25
$ws.Range("L128").Value = 57499.25
$ws.Range("N128").Value = -67459.25
$ws.Range("H132").Value = 1116.56
$ws.Range("I132").Value = 1116.56
$ws.Range("K132").Value = 3349.68
$ws.Range("M132").Value = -819.6799999999998
$ws.Range("H136").Value = 1480.9445
$ws.Range("I136").Value = 1416.125
$ws.Range("J136").Value = 1999.5
$ws.Range("K136").Value = 4248.375
$ws.Range("L136").Value = 5998.5
$ws.Range("M136").Value = -1698.375
$ws.Range("N136").Value = -11098.5
